$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string content updates -------------------------------

# A3: "红色格子：填入你的实验数据" -> "...，如本身自带数据请更改"
$ws.Range("A3").Value = "红色格子：填入你的实验数据，如本身自带数据请更改"

# A29: clarify that the table below is a plain helper table without formulas
$ws.Range("A29").Value = "过程中的温度（以下是个用于作图方便的纯表格，不包含公式, 选中后可以直接用excel散点作图）"

# A53: repo renamed from fuck-university-physics-experiments to fuck-nku-physics-experiments
$ws.Range("A53").Value = "Posted on https://github.com/Axolyz/fuck-nku-physics-experiments."

# A55: drop the trailing quote, leave the cell blank
$ws.Range("A55").ClearContents()

# --- Clear the placeholder "1 / 38 / 26" values in columns C & D --------
# (these were dummy numbers left in the template; the real template only
#  fills column B, leaving C/D blank for the user)

$ws.Range("C17:D19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("C23:D24").ClearContents()

# --- Update the remembered selection -------------------------------------
$null = $ws.Range("G19").Select()
